$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column C ("Förändrad") rows 2-264: update date serial value from 45172 to 45175
for ($r = 2; $r -le 264; $r++) {
    $cell = $ws.Cells.Item($r, 3)
    $v = $cell.Value2()
    if ($v -eq 45172) {
        $cell.Value = 45175
    }
}
